$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A172").Value = "06-09-2021"
$ws.Range("C172").Value = 4.47
